$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 364, shifting existing rows 364:390 down to 365:391.
$ws.Rows.Item(364).Insert()

# Populate the newly inserted row 364 with the new weekly data point.
$ws.Cells.Item(364, 1).Value  = 10
$ws.Cells.Item(364, 2).Value  = "Vega Modelo de Temuco"
$ws.Cells.Item(364, 3).Value  = "La Araucanía"
$ws.Cells.Item(364, 4).Value  = 44826
$ws.Cells.Item(364, 5).Value  = 9
$ws.Cells.Item(364, 6).Value  = 100112009
$ws.Cells.Item(364, 7).Value  = "Acelga"
$ws.Cells.Item(364, 8).Value  = "Sin especificar"
$ws.Cells.Item(364, 9).Value  = "Primera"
$ws.Cells.Item(364, 10).Value = 50
$ws.Cells.Item(364, 11).Value = 8000
$ws.Cells.Item(364, 12).Value = 8000
$ws.Cells.Item(364, 13).Value = 8000
$ws.Cells.Item(364, 14).Value = "$/docena de atados (12 kilos)"
$ws.Cells.Item(364, 15).Value = "Provincia de Cautín"
$ws.Cells.Item(364, 16).Value = 667
$ws.Cells.Item(364, 17).Value = 12
$ws.Cells.Item(364, 18).Value = "Hortaliza"
